# Read MACRO years directly from config of Excel
# Adds a "year" column (D) to the "config" sheet, listing the MACRO
# calibration years (2020, 2030, 2040) read from the model configuration.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")

# Header for the new column
$ws.Range("D1").Value = "year"

# MACRO years
$ws.Range("D2").Value = 2020
$ws.Range("D3").Value = 2030
$ws.Range("D4").Value = 2040

$ws.Activate()
$ws.Range("G14").Select()

$ws2 = $wb.Worksheets.Item("MERtoPPP")
$ws2.Activate()
$ws2.Range("L22").Select()

$ws.Activate()
